$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "hot fix. Wrong sample data attached." -- the header row had a stray
# "sex" label in column A (left over from bad sample data) and was
# missing a proper "Sex" header; fix the header labels and drop the
# bogus placeholder rows that shipped with the wrong sample export.

$ws.Range("A1").Value = "Date of Birth"
$ws.Range("B1").Value = "Date of Expiry"
$ws.Range("C1").Value = "Passport Number"
$ws.Range("D1").Value = "First Name"
$ws.Range("E1").Value = "Last Name"
$ws.Range("F1").Value = "Issuing State"
$ws.Range("G1").Value = "Nationality"
$ws.Range("H1").Value = "Sex"
$ws.Range("I1").Value = "NRIC Number"

# Remove the wrong sample data rows (2-10) entirely.
$ws.Rows("2:10").Delete()
